# Apply the "Context" sheet Landmark updates:
#  - split the lat/long JSON string (column G) into separate numeric
#    latitude (G) / longitude (H) columns for the existing ACES landmarks
#  - append six new Landmark rows (UMD / UNCC / CU campuses)
#  - bump the running total in A1
#  - restore the active-cell selection left by the editor

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Context")

# Running total of rows (was 30, now 36 after the six new landmarks)
$ws.Cells.Item(1, 1).Value = 36

# Existing ACES landmarks: replace the combined JSON text in column G with
# separate numeric latitude (G) / longitude (H) values.
$latlong = @(
    @(21, 39.195,    -106.821141),
    @(22, 39.195008,  -106.82138),
    @(23, 39.195428,  -106.822069),
    @(24, 39.195798,  -106.822166),
    @(25, 39.195998,  -106.821823),
    @(26, 39.19643,   -106.821629),
    @(27, 39.196613,  -106.82156),
    @(28, 39.196513,  -106.521265),
    @(29, 39.196368,  -106.821265),
    @(30, 39.195324,  -106.821227),
    @(31, 39.195216,  -106.821004)
)

foreach ($entry in $latlong) {
    $r = $entry[0]
    $lat = $entry[1]
    $lon = $entry[2]
    $ws.Cells.Item($r, 7).Value = $lat
    $ws.Cells.Item($r, 8).Value = $lon
}

# New Landmark rows for University of Maryland, UNC Charlotte and CU Boulder.
# Columns C (code) / D (title) / E (description) are written in the exact
# order the original author entered them so that newly-created shared
# strings land on the same indices as the target workbook.
function Set-LandmarkRow {
    param(
        [int]$row,
        [int]$num,
        [string]$code,
        [string]$title,
        [string]$desc,
        [string]$site,
        [double]$lat,
        [double]$lon,
        [string[]]$order
    )
    $ws.Cells.Item($row, 1).Value = $num
    $ws.Cells.Item($row, 2).Value = "Landmark"

    foreach ($col in $order) {
        if ($col -eq "C") { $ws.Cells.Item($row, 3).Value = $code }
        elseif ($col -eq "D") { $ws.Cells.Item($row, 4).Value = $title }
        elseif ($col -eq "E") { $ws.Cells.Item($row, 5).Value = $desc }
    }

    $ws.Cells.Item($row, 6).Value = $site
    $ws.Cells.Item($row, 7).Value = $lat
    $ws.Cells.Item($row, 8).Value = $lon
}

Set-LandmarkRow 32 31 "umd_landmark_hcil" "HCIL" "Human Computer Interaction Laboratory" "umd" 39.987901 -76.941599 @("C","D","E")
Set-LandmarkRow 33 32 "umd_landmark_avw" "AV Williams" "AV Williams, Computer Science, UMIACS" "umd" 38.990752 -76.936271 @("D","E","C")
Set-LandmarkRow 34 33 "umd_landmark_mall" "McKeldin Mall" "McKeldin Mall" "umd" 38.987134 -76.9403059 @("D","E","C")
Set-LandmarkRow 35 34 "uncc_landmark_woodward" "Woordward Hall" "Woodward Hall" "uncc" 35.3072387 -80.7353323 @("E","D","C")
Set-LandmarkRow 36 35 "cu_landmark_dlc" "Discovery Learning Center" "DLC" "cu" 40.007614 -105.261771 @("C","D","E")
Set-LandmarkRow 37 36 "cu_landmark_c4c" "Center for Community" "C4C - Center for Community" "cu" 40.004443 -105.26484 @("E","D","C")

# Leave the selection/active cell where the editor left it.
$ws.Activate()
$ws.Range("H32").Select()
